$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Repair Kit" description (row 13, column D) to mention repairing Destroyed armour
$ws.Range("D13").Value = "A set of tools for performing minor repairs to armour and clothing. Proficiency allows you to halve the time required to repair a set to full strength, and allows you to repair {\\it Destroyed} armour. "

# Move the selection/active cell to D2 (matches the saved view state in the edit)
$ws.Range("D2").Select() | Out-Null
